$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44806
$ws.Range("K2").Value2 = 7000
$ws.Range("L2").Value2 = 7500
$ws.Range("M2").Value2 = 7250
$ws.Range("O2").Value2 = 'Provincia de Diguillín'
$ws.Range("P2").Value2 = 725
$ws.Range("D3").Value2 = 44819
$ws.Range("J3").Value2 = 100
$ws.Range("M3").Value2 = 7500
$ws.Range("O3").Value2 = 'Provincia de Diguillín'
$ws.Range("P3").Value2 = 750
$ws.Range("D4").Value2 = 44211
$ws.Range("J4").Value2 = 28
$ws.Range("K4").Value2 = 8000
$ws.Range("L4").Value2 = 8500
$ws.Range("M4").Value2 = 8214
$ws.Range("O4").Value2 = 'Región Metropolitana'
$ws.Range("P4").Value2 = 821
$ws.Range("D5").Value2 = 44838
$ws.Range("J5").Value2 = 120
$ws.Range("K5").Value2 = 6500
$ws.Range("M5").Value2 = 6750
$ws.Range("P5").Value2 = 675
$ws.Range("D6").Value2 = 44798
$ws.Range("L6").Value2 = 7000
$ws.Range("M6").Value2 = 7000
$ws.Range("P6").Value2 = 700
$ws.Range("D7").Value2 = 44790
$ws.Range("J7").Value2 = 60
$ws.Range("K7").Value2 = 8500
$ws.Range("L7").Value2 = 9000
$ws.Range("M7").Value2 = 8750
$ws.Range("O7").Value2 = 'Región Metropolitana'
$ws.Range("P7").Value2 = 875
$ws.Range("D8").Value2 = 44980
$ws.Range("K8").Value2 = 7500
$ws.Range("L8").Value2 = 8000
$ws.Range("M8").Value2 = 7750
$ws.Range("P8").Value2 = 775
$ws.Range("D9").Value2 = 44784
$ws.Range("J9").Value2 = 100
$ws.Range("K9").Value2 = 8000
$ws.Range("L9").Value2 = 9000
$ws.Range("M9").Value2 = 8500
$ws.Range("O9").Value2 = 'Región Metropolitana'
$ws.Range("P9").Value2 = 850
$ws.Range("D10").Value2 = 44782
$ws.Range("J10").Value2 = 120
$ws.Range("D11").Value2 = 44847
$ws.Range("J11").Value2 = 100
$ws.Range("K11").Value2 = 6500
$ws.Range("L11").Value2 = 7000
$ws.Range("M11").Value2 = 6750
$ws.Range("O11").Value2 = 'Provincia de Diguillín'
$ws.Range("P11").Value2 = 675
$ws.Range("D12").Value2 = 44804
$ws.Range("J12").Value2 = 80
$ws.Range("K12").Value2 = 7000
$ws.Range("L12").Value2 = 7500
$ws.Range("M12").Value2 = 7250
$ws.Range("P12").Value2 = 725
$ws.Range("D13").Value2 = 44810
$ws.Range("J13").Value2 = 60
$ws.Range("K13").Value2 = 7000
$ws.Range("L13").Value2 = 8000
$ws.Range("M13").Value2 = 7500
$ws.Range("P13").Value2 = 750
$ws.Range("D14").Value2 = 44791
$ws.Range("J14").Value2 = 100
$ws.Range("K14").Value2 = 8500
$ws.Range("L14").Value2 = 9000
$ws.Range("M14").Value2 = 8750
$ws.Range("O14").Value2 = 'Región Metropolitana'
$ws.Range("P14").Value2 = 875
$ws.Range("D15").Value2 = 44775
$ws.Range("J15").Value2 = 60
$ws.Range("K15").Value2 = 8000
$ws.Range("L15").Value2 = 8000
$ws.Range("M15").Value2 = 8000
$ws.Range("O15").Value2 = 'Región Metropolitana'
$ws.Range("P15").Value2 = 800
$ws.Range("D16").Value2 = 44817
$ws.Range("K16").Value2 = 7000
$ws.Range("L16").Value2 = 7000
$ws.Range("M16").Value2 = 7000
$ws.Range("O16").Value2 = 'Provincia de Diguillín'
$ws.Range("P16").Value2 = 700
$ws.Range("D17").Value2 = 44817
$ws.Range("I17").Value2 = 'Segunda'
$ws.Range("J17").Value2 = 60
$ws.Range("K17").Value2 = 8000
$ws.Range("M17").Value2 = 8000
$ws.Range("P17").Value2 = 800
$ws.Range("D18").Value2 = 44831
$ws.Range("J18").Value2 = 60
$ws.Range("K18").Value2 = 7000
$ws.Range("L18").Value2 = 7500
$ws.Range("M18").Value2 = 7250
$ws.Range("O18").Value2 = 'Provincia de Diguillín'
$ws.Range("P18").Value2 = 725
$ws.Range("D19").Value2 = 44799
$ws.Range("L19").Value2 = 7000
$ws.Range("M19").Value2 = 7000
$ws.Range("P19").Value2 = 700
$ws.Range("D20").Value2 = 44841
$ws.Range("K20").Value2 = 6500
$ws.Range("L20").Value2 = 7000
$ws.Range("M20").Value2 = 6750
$ws.Range("P20").Value2 = 675
$ws.Range("D21").Value2 = 44813
$ws.Range("J21").Value2 = 120
$ws.Range("K21").Value2 = 7000
$ws.Range("L21").Value2 = 7500
$ws.Range("M21").Value2 = 7250
$ws.Range("O21").Value2 = 'Provincia de Diguillín'
$ws.Range("P21").Value2 = 725
$ws.Range("D22").Value2 = 44846
$ws.Range("J22").Value2 = 100
$ws.Range("K22").Value2 = 6500
$ws.Range("L22").Value2 = 7000
$ws.Range("M22").Value2 = 6750
$ws.Range("P22").Value2 = 675
$ws.Range("D23").Value2 = 44203
$ws.Range("J23").Value2 = 27
$ws.Range("L23").Value2 = 8000
$ws.Range("M23").Value2 = 7556
$ws.Range("O23").Value2 = 'Región Metropolitana'
$ws.Range("P23").Value2 = 756
$ws.Range("D24").Value2 = 44812
$ws.Range("I24").Value2 = 'Primera'
$ws.Range("K24").Value2 = 7000
$ws.Range("M24").Value2 = 7500
$ws.Range("P24").Value2 = 750
